# 10Th - MB for single stock and added new group
#
# This workbook is a "MarketBeat rank" watch sheet: column A holds analyst
# names, and each subsequent column holds that analyst's most-recent rating
# action for one historical date (newest date on the left in column B,
# older dates to the right). This edit adds two brand-new date columns
# (two new "as of" snapshots) in front of the existing date columns, and
# appends two new analyst rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out the current extent of the table before we start mutating it.
$used      = $ws.UsedRange
$lastRow   = $used.Rows.Count          # 27 -> rows 1..27 already populated
$firstDataRow = 2                      # row 1 is the header row

# --- 1. Make room for three new "date" columns -----------------------
# The existing B:E columns (one per historical date) all shift right by
# three columns, becoming E:H, carrying their values/styles with them.
$ws.Range("B:D").Insert()

# --- 2. New header row values ------------------------------------------
# B1/C1/D1 are the three brand-new date columns.
$ws.Cells.Item(1, 2).Value = "Jun_27"
$ws.Cells.Item(1, 3).Value = "Jun_26"
$ws.Cells.Item(1, 4).Value = "Jun_26"

# --- 3. Fill the new columns for every existing analyst row with the
#        default "UN" (unchanged) marker, same as every other empty cell
#        in the grid.
for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# --- 4. Add the new analyst group at the bottom of the table -----------
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

$ws.Cells.Item($newRow1, 1).Value = "Benchmark"
$ws.Cells.Item($newRow1, 2).Value = "UN"
$ws.Cells.Item($newRow1, 3).Value = "UN"
$ws.Cells.Item($newRow1, 4).Value = "UN"

$ws.Cells.Item($newRow2, 1).Value = "Evercore ISI"
$ws.Cells.Item($newRow2, 2).Value = "UN"
$ws.Cells.Item($newRow2, 3).Value = "UN"
$ws.Cells.Item($newRow2, 4).Value = "UN"
